$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CH4")
$ws1.Rows(15).Copy()
$ws1.Rows(16).PasteSpecial(-4122)
Write-Host "done"
